$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy label style (bold/border/center) from an existing label cell (A2) for the 5 newly appended rows (27-31).
$ws.Range("A2").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)

# Row 17: NB_F
$ws.Range("A17").Value = "NB_F"
$ws.Range("B17").Value = 0.8775025799793603
$ws.Range("C17").Value = 0.9183006535947712
$ws.Range("D17").Value = 0.8970424205718324
$ws.Range("E17").Value = 0.878809769521844
$ws.Range("F17").Value = 0.9300653594771241
$ws.Range("G17").Value = 0.903428975193681
$ws.Range("H17").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()
$ws.Range("K17").Value = 0.8669762641898865
$ws.Range("L17").Value = 0.9065359477124183
$ws.Range("M17").Value = 0.8859313094607213
$ws.Range("N17").Value = 0.8795846233230135
$ws.Range("O17").Value = 0.9183006535947712
$ws.Range("P17").Value = 0.8972849090496149
$ws.Range("Q17").ClearContents()
$ws.Range("R17").ClearContents()
$ws.Range("S17").ClearContents()

# Row 18: NB_NF
$ws.Range("A18").Value = "NB_NF"
$ws.Range("B18").Value = 0.9421826625386996
$ws.Range("C18").Value = 0.8078947368421051
$ws.Range("D18").Value = 0.8692109692109693
$ws.Range("E18").Value = 0.9546826625386997
$ws.Range("F18").Value = 0.8078947368421051
$ws.Range("G18").Value = 0.873972873972874
$ws.Range("H18").ClearContents()
$ws.Range("I18").ClearContents()
$ws.Range("J18").ClearContents()
$ws.Range("K18").Value = 0.9311532507739937
$ws.Range("L18").Value = 0.7978947368421052
$ws.Range("M18").Value = 0.8587004587004585
$ws.Range("N18").Value = 0.9311532507739937
$ws.Range("O18").Value = 0.7978947368421052
$ws.Range("P18").Value = 0.8587004587004585
$ws.Range("Q18").ClearContents()
$ws.Range("R18").ClearContents()
$ws.Range("S18").ClearContents()

# Row 19: NB_PM
$ws.Range("A19").Value = "NB_PM"
$ws.Range("B19").Value = 0.7925077211118401
$ws.Range("C19").Value = 0.9263157894736842
$ws.Range("D19").Value = 0.8520716425209365
$ws.Range("E19").Value = 0.799064935064935
$ws.Range("F19").Value = 0.9368421052631579
$ws.Range("G19").Value = 0.8610469545580154
$ws.Range("H19").ClearContents()
$ws.Range("I19").ClearContents()
$ws.Range("J19").ClearContents()
$ws.Range("K19").Value = 0.7386363636363636
$ws.Range("L19").Value = 0.9368421052631579
$ws.Range("M19").Value = 0.8251547510419492
$ws.Range("N19").Value = 0.7534914361001317
$ws.Range("O19").Value = 0.9473684210526315
$ws.Range("P19").Value = 0.8384437875867997
$ws.Range("Q19").ClearContents()
$ws.Range("R19").ClearContents()
$ws.Range("S19").ClearContents()

# Row 20: NB_M
$ws.Range("A20").Value = "NB_M"
$ws.Range("B20").Value = 0.825701754385965
$ws.Range("C20").Value = 0.7485294117647059
$ws.Range("D20").Value = 0.781008921248076
$ws.Range("E20").Value = 0.8380392156862746
$ws.Range("F20").Value = 0.7838235294117647
$ws.Range("G20").Value = 0.8094028520499108
$ws.Range("H20").ClearContents()
$ws.Range("I20").ClearContents()
$ws.Range("J20").ClearContents()
$ws.Range("K20").Value = 0.8564102564102564
$ws.Range("L20").Value = 0.6294117647058822
$ws.Range("M20").Value = 0.7201566555431961
$ws.Range("N20").Value = 0.8597435897435897
$ws.Range("O20").Value = 0.6529411764705881
$ws.Range("P20").Value = 0.7385474601408972
$ws.Range("Q20").ClearContents()
$ws.Range("R20").ClearContents()
$ws.Range("S20").ClearContents()

# Row 21: NB_V
$ws.Range("A21").Value = "NB_V"
$ws.Range("B21").Value = 0.8347924158221641
$ws.Range("C21").Value = 0.8342857142857143
$ws.Range("D21").Value = 0.8332294062526622
$ws.Range("E21").Value = 0.8609394698868382
$ws.Range("F21").Value = 0.8438095238095238
$ws.Range("G21").Value = 0.8513013105651253
$ws.Range("H21").ClearContents()
$ws.Range("I21").ClearContents()
$ws.Range("J21").ClearContents()
$ws.Range("K21").Value = 0.8088104153664795
$ws.Range("L21").Value = 0.8538095238095238
$ws.Range("M21").Value = 0.8287315010570824
$ws.Range("N21").Value = 0.823794466403162
$ws.Range("O21").Value = 0.8633333333333333
$ws.Range("P21").Value = 0.8421883693919995
$ws.Range("Q21").ClearContents()
$ws.Range("R21").ClearContents()
$ws.Range("S21").ClearContents()

# Row 22: RF_F
$ws.Range("A22").Value = "RF_F"
$ws.Range("B22").Value = 0.8856372549019607
$ws.Range("C22").Value = 0.8026143790849674
$ws.Range("D22").Value = 0.8412433155080216
$ws.Range("E22").Value = 0.8617927170868347
$ws.Range("F22").Value = 0.7797385620915032
$ws.Range("G22").Value = 0.8169117647058822
$ws.Range("H22").Value = 0.8375101214574899
$ws.Range("I22").Value = 0.8124183006535948
$ws.Range("J22").Value = 0.8207655502392344
$ws.Range("K22").Value = 0.8714705882352941
$ws.Range("L22").Value = 0.7908496732026145
$ws.Range("M22").Value = 0.8283645276292335
$ws.Range("N22").Value = 0.8766666666666666
$ws.Range("O22").Value = 0.8143790849673203
$ws.Range("P22").Value = 0.8430194805194805
$ws.Range("Q22").Value = 0.870031696889282
$ws.Range("R22").Value = 0.9065359477124183
$ws.Range("S22").Value = 0.886539604186663

# Row 23: RF_NF
$ws.Range("A23").Value = "RF_NF"
$ws.Range("B23").Value = 0.7099472990777338
$ws.Range("C23").Value = 0.8384210526315791
$ws.Range("D23").Value = 0.768448605657908
$ws.Range("E23").Value = 0.7424987864919215
$ws.Range("F23").Value = 0.8278947368421055
$ws.Range("G23").Value = 0.7813279743512301
$ws.Range("H23").Value = 0.7742181540808544
$ws.Range("I23").Value = 0.8173684210526316
$ws.Range("J23").Value = 0.7911215915867078
$ws.Range("K23").Value = 0.7301739130434782
$ws.Range("L23").Value = 0.8384210526315791
$ws.Range("M23").Value = 0.7791754756871037
$ws.Range("N23").Value = 0.7081501831501831
$ws.Range("O23").Value = 0.848421052631579
$ws.Range("P23").Value = 0.7709426754162753
$ws.Range("Q23").Value = 0.8454720133667502
$ws.Range("R23").Value = 0.8278947368421052
$ws.Range("S23").Value = 0.8360156018564234

# Row 24: RF_PM
$ws.Range("A24").Value = "RF_PM"
$ws.Range("B24").Value = 0.8541197579510271
$ws.Range("C24").Value = 0.8315789473684211
$ws.Range("D24").Value = 0.8402062065733439
$ws.Range("E24").Value = 0.8480247678018575
$ws.Range("F24").Value = 0.8736842105263157
$ws.Range("G24").Value = 0.857203615098352
$ws.Range("H24").Value = 0.7741993888464477
$ws.Range("I24").Value = 0.8105263157894737
$ws.Range("J24").Value = 0.7866198571076619
$ws.Range("K24").Value = 0.8303581480439066
$ws.Range("L24").Value = 0.8105263157894737
$ws.Range("M24").Value = 0.8179338590378386
$ws.Range("N24").Value = 0.8050597508263412
$ws.Range("O24").Value = 0.8210526315789475
$ws.Range("P24").Value = 0.8112578744157692
$ws.Range("Q24").Value = 0.7748917748917749
$ws.Range("R24").Value = 0.8736842105263157
$ws.Range("S24").Value = 0.821219512195122

# Row 25: RF_M
$ws.Range("A25").Value = "RF_M"
$ws.Range("B25").Value = 0.7616410912190965
$ws.Range("C25").Value = 0.8102941176470588
$ws.Range("D25").Value = 0.7793885102708632
$ws.Range("E25").Value = 0.7439164086687307
$ws.Range("F25").Value = 0.7624999999999998
$ws.Range("G25").Value = 0.751029836176895
$ws.Range("H25").Value = 0.7496031746031747
$ws.Range("I25").Value = 0.6786764705882352
$ws.Range("J25").Value = 0.7088569394791996
$ws.Range("K25").Value = 0.731552250190694
$ws.Range("L25").Value = 0.7977941176470589
$ws.Range("M25").Value = 0.7579477204477205
$ws.Range("N25").Value = 0.7645704948646126
$ws.Range("O25").Value = 0.7625
$ws.Range("P25").Value = 0.7598665429025961
$ws.Range("Q25").Value = 0.7678991596638656
$ws.Range("R25").Value = 0.7727941176470587
$ws.Range("S25").Value = 0.7668180785827845

# Row 26: RF_V
$ws.Range("A26").Value = "RF_V"
$ws.Range("B26").Value = 0.9298331613347093
$ws.Range("C26").Value = 0.7947619047619047
$ws.Range("D26").Value = 0.8559604138551506
$ws.Range("E26").Value = 0.9041083099906629
$ws.Range("F26").Value = 0.8147619047619047
$ws.Range("G26").Value = 0.8555017938843354
$ws.Range("H26").Value = 0.7802756892230577
$ws.Range("I26").Value = 0.7557142857142857
$ws.Range("J26").Value = 0.7670933431138824
$ws.Range("K26").Value = 0.9111111111111111
$ws.Range("L26").Value = 0.7842857142857144
$ws.Range("M26").Value = 0.840944669365722
$ws.Range("N26").Value = 0.9075421396628827
$ws.Range("O26").Value = 0.7547619047619047
$ws.Range("P26").Value = 0.8222915222915222
$ws.Range("Q26").Value = 0.8629411764705882
$ws.Range("R26").Value = 0.7271428571428571
$ws.Range("S26").Value = 0.7875486012328119

# Row 27: Ensemble_F
$ws.Range("A27").Value = "Ensemble_F"
$ws.Range("B27").Value = 0.8775025799793603
$ws.Range("C27").Value = 0.9183006535947712
$ws.Range("D27").Value = 0.8970424205718324
$ws.Range("E27").Value = 0.8883006535947711
$ws.Range("F27").Value = 0.9183006535947712
$ws.Range("G27").Value = 0.902084437378555
$ws.Range("H27").Value = 0.8192676767676769
$ws.Range("I27").Value = 0.8124183006535948
$ws.Range("J27").Value = 0.8125757575757575
$ws.Range("K27").Value = 0.8669762641898865
$ws.Range("L27").Value = 0.9065359477124183
$ws.Range("M27").Value = 0.8859313094607213
$ws.Range("N27").Value = 0.8346491228070176
$ws.Range("O27").Value = 0.9294117647058823
$ws.Range("P27").Value = 0.8777041671778514
$ws.Range("Q27").Value = 0.8572222222222223
$ws.Range("R27").Value = 0.8477124183006535
$ws.Range("S27").Value = 0.8514430014430016

# Row 28: Ensemble_NF
$ws.Range("A28").Value = "Ensemble_NF"
$ws.Range("B28").Value = 0.9421826625386996
$ws.Range("C28").Value = 0.8078947368421051
$ws.Range("D28").Value = 0.8692109692109693
$ws.Range("E28").Value = 0.943888888888889
$ws.Range("F28").Value = 0.8178947368421052
$ws.Range("G28").Value = 0.8748203842940685
$ws.Range("H28").Value = 0.7932900432900433
$ws.Range("I28").Value = 0.8263157894736842
$ws.Range("J28").Value = 0.8075261324041814
$ws.Range("K28").Value = 0.9311532507739937
$ws.Range("L28").Value = 0.7978947368421052
$ws.Range("M28").Value = 0.8587004587004585
$ws.Range("N28").Value = 0.9307189542483659
$ws.Range("O28").Value = 0.758421052631579
$ws.Range("P28").Value = 0.8336846355422207
$ws.Range("Q28").Value = 0.8142395762132605
$ws.Range("R28").Value = 0.7978947368421052
$ws.Range("S28").Value = 0.8029768605378361

# Row 29: Ensemble_PM
$ws.Range("A29").Value = "Ensemble_PM"
$ws.Range("B29").Value = 0.7925077211118401
$ws.Range("C29").Value = 0.9263157894736842
$ws.Range("D29").Value = 0.8520716425209365
$ws.Range("E29").Value = 0.8119138755980861
$ws.Range("F29").Value = 0.8842105263157893
$ws.Range("G29").Value = 0.8447665041568972
$ws.Range("H29").Value = 0.7786354775828459
$ws.Range("I29").Value = 0.8526315789473685
$ws.Range("J29").Value = 0.809889674924
$ws.Range("K29").Value = 0.7386363636363636
$ws.Range("L29").Value = 0.9368421052631579
$ws.Range("M29").Value = 0.8251547510419492
$ws.Range("N29").Value = 0.785972360972361
$ws.Range("O29").Value = 0.9473684210526315
$ws.Range("P29").Value = 0.8579636982416335
$ws.Range("Q29").Value = 0.8157302346776032
$ws.Range("R29").Value = 0.8631578947368421
$ws.Range("S29").Value = 0.8380948619202793

# Row 30: Ensemble_M
$ws.Range("A30").Value = "Ensemble_M"
$ws.Range("B30").Value = 0.825701754385965
$ws.Range("C30").Value = 0.7485294117647059
$ws.Range("D30").Value = 0.781008921248076
$ws.Range("E30").Value = 0.7983900928792569
$ws.Range("F30").Value = 0.7970588235294118
$ws.Range("G30").Value = 0.7949042950513538
$ws.Range("H30").Value = 0.8066793681035167
$ws.Range("I30").Value = 0.6904411764705882
$ws.Range("J30").Value = 0.7344887955182072
$ws.Range("K30").Value = 0.8564102564102564
$ws.Range("L30").Value = 0.6294117647058822
$ws.Range("M30").Value = 0.7201566555431961
$ws.Range("N30").Value = 0.910989010989011
$ws.Range("O30").Value = 0.7139705882352941
$ws.Range("P30").Value = 0.7999999999999999
$ws.Range("Q30").Value = 0.8200962436256554
$ws.Range("R30").Value = 0.7132352941176471
$ws.Range("S30").Value = 0.7587602783048704

# Row 31: Ensemble_V
$ws.Range("A31").Value = "Ensemble_V"
$ws.Range("B31").Value = 0.8347924158221641
$ws.Range("C31").Value = 0.8342857142857143
$ws.Range("D31").Value = 0.8332294062526622
$ws.Range("E31").Value = 0.8639553429027114
$ws.Range("F31").Value = 0.8638095238095238
$ws.Range("G31").Value = 0.8628499955329223
$ws.Range("H31").Value = 0.7757787050747071
$ws.Range("I31").Value = 0.7452380952380953
$ws.Range("J31").Value = 0.7579824010056568
$ws.Range("K31").Value = 0.8088104153664795
$ws.Range("L31").Value = 0.8538095238095238
$ws.Range("M31").Value = 0.8287315010570824
$ws.Range("N31").Value = 0.8558095238095238
$ws.Range("O31").Value = 0.9023809523809524
$ws.Range("P31").Value = 0.8776215724890168
$ws.Range("Q31").Value = 0.8062280701754385
$ws.Range("R31").Value = 0.8442857142857143
$ws.Range("S31").Value = 0.8219733924611974
